$oldTitleHe = "המאמר היומי של מייק - 04.04.25"
$newTitleHe = "המאמר היומי של מייק - 02.04.25"
$oldTitleEn = "Amortizing intractable inference in diffusion models for vision, language, and control"
$newTitleEn = "SymDPO: Boosting In-Context Learning of Large Multimodal Models with Symbol Demonstration Direct Preference Optimization"
$newP2 = "היום אני עושה מעבר חד בנושא הסקירה וסוקר מאמר על אימון מודלים מולטימודליים (בפרט MLLMs). המאמר מציע שיטה לאימון מודלים למשימת למידה in-context שבא המודל מקבל כמה דוגמאות (הדגמות) שכל הדגמה מכילה תמונה, שאלה ותשובה עליה. המודל מתבקש, בהתבסס על הדאטה שקיבל (הדגמות) לענות על שאלה לגבי תמונה נוספת (עם אותם הדמויות למשל). הסקירה הולכת להיות קלילה וקצרה."
$newP3 = "המחברים מציעים דרך לשיפור הבנת קשרים בין פיסות דאטה ממודליות שונות על ידי מודלים מולטי-מודליים. למשל למודלים התומכים בשתי מודליות, שפה ותמונות, לפעמים מתקשים במשימות שדורשות הבנת קשרים סמנטיים בין דאטה ויזואלי לשפתי למשל במשימת למידה in-context ל-MLLMs המתוארת קודם לכם. המאמר מציין כי MLLMs לפעמים מתקשים להתמודד עם משימות אלו ולמשל עונים על השאלה בלי להתחשב בהקשר כלל (שזה תמונות, שאלות ותשובות). המאמר מציע שיטת פיין טיון עבור מודל מולטי-מודלי כדי להתמודד עם כשלים כאלו."
$newP4 = "המאמר מציע לעשות פיינטיון למודל בשיטה מעולם RLHF (שזה Reinforcement Learning with Human Feedback) הנקראת(DPO (= Direct Preference Optimization . שיטה זו נגזרת מפונקציית יעד פופולרית בעבור פיין טיון של מודל שפה  (מקסום תגמול - קרבה למודל ההתחלתי) דרך מידול reward של Bradley-Terry. היתרון העיקרי של DPO מעל PPO הוא העובדה ש-DPO לא דורש אימון של מודל תגמול (reward) אלא צריך רק דאטהסט של זוגות שאלות ותשובות רצויות ותשובות לא רצויות. הרעיון העיקרי במאמר הוא להנדס דאטהסט כזה עבור יוזקייס שבנידון ולהשתמש ב-DPO לפיין טיון של מודל מולטימודלי."
$newP5 = "בגדול המאמר מציע לשחק עם השאלות והתשובות. הוא מציע כמה טריקים כדי לאלץ את המודלי להתחשב בכל הקונטקסט שניתן לו. אחד הטריקים הוא לתת תשובה רצויה לא קשורה (מילה ללא משמעות). עוד טריק היא להחליף תשובה לא רצויה בג'בריש ועוד אחד היא למחוק את השאלה עצמה ולהשאיר את התשובות כמו שהם. יש עוד כמה טריקים מהסוג הזה ועל ידי שילובם המאמר משיג מודל יותר טוב עם שימוש ב-DPO לפיין-טיון. "
$newP6 = "כמו שהבטחתי סקירה קצרה וקלילה."
$newUrl = "https://arxiv.org/abs/2411.11909"

$d = $word.ActiveDocument

# --- Paragraph 1 (title paragraph): two <w:t> runs separated by a line break ---
$d.Content.Find.Execute($oldTitleHe, $true, $false, $false, $false, $false, $true, 1, $false, $newTitleHe, 2) | Out-Null
$d.Content.Find.Execute($oldTitleEn, $true, $false, $false, $false, $false, $true, 1, $false, $newTitleEn, 2) | Out-Null

# --- Paragraphs 2-6: full paragraph text replacement (exclude trailing paragraph mark) ---
$p2 = $d.Paragraphs.Item(2)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$r2.Text = $newP2

$p3 = $d.Paragraphs.Item(3)
$r3 = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$r3.Text = $newP3

$p4 = $d.Paragraphs.Item(4)
$r4 = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$r4.Text = $newP4

$p5 = $d.Paragraphs.Item(5)
$r5 = $d.Range($p5.Range.Start, $p5.Range.End - 1)
$r5.Text = $newP5

$p6 = $d.Paragraphs.Item(6)
$r6 = $d.Range($p6.Range.Start, $p6.Range.End - 1)
$r6.Text = $newP6

# --- Delete paragraphs 7 through 31 (the old deep-dive body + old URL paragraph) ---
$startPara = $d.Paragraphs.Item(7)
$endPara = $d.Paragraphs.Item(31)
$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()

# --- Remaining last paragraph is now the URL paragraph; update its link text ---
$lastIndex = $d.Paragraphs.Count
$pUrl = $d.Paragraphs.Item($lastIndex)
$rUrl = $d.Range($pUrl.Range.Start, $pUrl.Range.End - 1)
$rUrl.Text = $newUrl

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
